# Update the Power Supply spec from 24V, 3Amp to 24V, 5Amp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "24V, 5Amp"

# Update the selected/active cell to B22 (matches the saved view state in the diff)
$ws.Range("B22").Select()
